$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (244) down into
# the three new rows so the new date cells in column A keep the same style
# (s="2", date number format) as the rest of the column.
$ws.Range("A244").Copy($ws.Range("A245:A247"))

# Row 245 -> 2021-05-03 (serial 44319)
$ws.Range("A245").Value = 44319
$ws.Range("B245").Value = 0
$ws.Range("C245").Value = 0
$ws.Range("D245").Value = 0

# Row 246 -> 2021-05-04 (serial 44320)
$ws.Range("A246").Value = 44320
$ws.Range("B246").Value = 0
$ws.Range("C246").Value = 0
$ws.Range("D246").Value = 0

# Row 247 -> 2021-05-05 (serial 44321)
$ws.Range("A247").Value = 44321
$ws.Range("B247").Value = 0
$ws.Range("C247").Value = 0
$ws.Range("D247").Value = 0
